# Weekly fruit/vegetable price update: a new price observation was
# recorded, so insert a new row at row 21 (pushing the existing rows
# 21-41 down to 22-42) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("21:21").Insert()

$ws.Cells.Item(21, 1).Value2  = 10
$ws.Cells.Item(21, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value2  = "La Araucanía"
$ws.Cells.Item(21, 4).Value2  = 44587
$ws.Cells.Item(21, 5).Value2  = 9
$ws.Cells.Item(21, 6).Value2  = 100114002
$ws.Cells.Item(21, 7).Value2  = "Camote"
$ws.Cells.Item(21, 8).Value2  = "Sin especificar"
$ws.Cells.Item(21, 9).Value2  = "Primera"
$ws.Cells.Item(21, 10).Value2 = 55
$ws.Cells.Item(21, 11).Value2 = 18000
$ws.Cells.Item(21, 12).Value2 = 18000
$ws.Cells.Item(21, 13).Value2 = 18000
$ws.Cells.Item(21, 14).Value2 = "`$/malla 20 kilos"
$ws.Cells.Item(21, 15).Value2 = "Perú"
$ws.Cells.Item(21, 16).Value2 = 900
$ws.Cells.Item(21, 17).Value2 = 20
$ws.Cells.Item(21, 18).Value2 = "Hortaliza"
